$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, preserving them as text
$priceUpdates = @{
    "D2" = "64.711.75"
    "D3" = "3.148.67"
    "D5" = "576.62"
    "D6" = "148.60"
    "D8" = "3.145.64"
    "D14" = "37.06"
    "D15" = "3.668.43"
    "D16" = "64.829.23"
    "D17" = "3.176.55"
    "D20" = "503.54"
    "D22" = "15.30"
    "D24" = "7.69"
    "D25" = "83.91"
    "D26" = "1.00"
    "D27" = "8.89"
    "D31" = "27.42"
    "D32" = "0.999"
    "D34" = "6.17"
    "D36" = "54.59"
    "D38" = "476.02"
    "D41" = "8.67"
    "D42" = "2.990.82"
    "D46" = "28.04"
    "D47" = "0.0₃0580"
    "D51" = "33.50"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Update Volume(1h) (column E) percentage values
$volumeUpdates = @{
    "E2" = "  -0.21%  "
    "E3" = "  +0.06%  "
    "E4" = "  -0.03%  "
    "E5" = "  +0.76%  "
    "E6" = "  -1.70%  "
    "E7" = "  +0.13%  "
    "E8" = "  +0.00%  "
    "E9" = "  -0.45%  "
    "E10" = "  -2.00%  "
    "E11" = "  -1.10%  "
    "E12" = "  -0.82%  "
    "E13" = "  +2.51%  "
    "E14" = "  -1.27%  "
    "E15" = "  +0.25%  "
    "E16" = "  -0.15%  "
    "E17" = "  +0.90%  "
    "E18" = "  -1.79%  "
    "E19" = "  +0.35%  "
    "E20" = "  -1.58%  "
    "E21" = "  -0.63%  "
    "E22" = "  -0.18%  "
    "E23" = "  -3.11%  "
    "E24" = "  -1.58%  "
    "E26" = "  +0.11%  "
    "E27" = "  +2.02%  "
    "E28" = "  -1.28%  "
    "E29" = "  -1.08%  "
    "E30" = "  +5.89%  "
    "E31" = "  -2.00%  "
    "E32" = "  -0.08%  "
    "E33" = "  +1.17%  "
    "E34" = "  +1.32%  "
    "E35" = "  -1.98%  "
    "E36" = "  -1.52%  "
    "E37" = "  +3.33%  "
    "E38" = "  -1.37%  "
    "E39" = "  -2.22%  "
    "E40" = "  -2.78%  "
    "E41" = "  +0.34%  "
    "E42" = "  -4.00%  "
    "E43" = "  -3.88%  "
    "E44" = "  -2.11%  "
    "E45" = "  -4.07%  "
    "E46" = "  -3.61%  "
    "E47" = "  +1.28%  "
    "E49" = "  -1.65%  "
    "E50" = "  -2.82%  "
    "E51" = "  +6.77%  "
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

